# 21/04/2023: dB: Login se configura con dB y queda funcionando
#
# 1) "Labels" sheet: insert a new "Id" column at the front, numbering the
#    existing label rows 1..9.
# 2) "Database" sheet: update the dBUser / dBPass values with the working
#    credentials, and make it the active sheet/cell.
# 3) "Usuarios RMobile" sheet: just a leftover cell-selection change.
# 4) "DbTables" sheet stops being the active sheet (Database becomes active
#    instead) - this happens automatically once Database is activated.

$wb = $excel.ActiveWorkbook

# --- Labels sheet: insert Id column -----------------------------------
$wsLabels = $wb.Worksheets.Item("Labels")
$null = $wsLabels.Columns.Item(1).Insert()
$wsLabels.Range("A1").Value = "Id"
for ($i = 1; $i -le 9; $i++) {
    $wsLabels.Cells.Item($i + 1, 1).Value = $i
}
$null = $wsLabels.Range("A11").Select()

# --- Usuarios RMobile sheet: selection only ----------------------------
$wsUsuarios = $wb.Worksheets.Item("Usuarios RMobile")
$null = $wsUsuarios.Range("F17").Select()

# --- Database sheet: update credentials + become active sheet ----------
$wsDatabase = $wb.Worksheets.Item("Database")
$wsDatabase.Range("B4").Value = "k4v441pbw49j1bqk4gbl"
$wsDatabase.Range("B5").Value = "pscale_pw_mYGOaWhbbqWTnBjKQo5zhijPujvef0XU7ggh5eQ0tCC"
$null = $wsDatabase.Activate()
$null = $wsDatabase.Range("B5").Select()
